$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row (row 1) had columns CR1:DN1 holding truncated/placeholder
# category labels (DESCRIÇÃO ESTADO CIVIL_<first letter>, NACIONALIDADE_<code>).
# Replace them with the full, correct category labels. The corrected list is
# shorter (21 columns instead of 23), so the two trailing columns are cleared.

$newHeaders = @(
    "DESCRIÇÃO ESTADO CIVIL_Casado",
    "DESCRIÇÃO ESTADO CIVIL_Desquitado",
    "DESCRIÇÃO ESTADO CIVIL_Divorciado",
    "DESCRIÇÃO ESTADO CIVIL_Outros",
    "DESCRIÇÃO ESTADO CIVIL_Separado",
    "DESCRIÇÃO ESTADO CIVIL_Solteiro",
    "DESCRIÇÃO ESTADO CIVIL_União Estável",
    "DESCRIÇÃO ESTADO CIVIL_Viúvo",
    "NACIONALIDADE_Brasileiro",
    "NACIONALIDADE_Estrangeiro",
    "GRAUINSTRUCAO_5º ano completo do ensino fundamental",
    "GRAUINSTRUCAO_Analfabeto",
    "GRAUINSTRUCAO_Até o 5º ano incompleto do ensino fundamental",
    "GRAUINSTRUCAO_Do 6º ao 9º ano do ensino fundamental ",
    "GRAUINSTRUCAO_Educação superior completo",
    "GRAUINSTRUCAO_Educação superior incompleto",
    "GRAUINSTRUCAO_Ensino fundamental completo",
    "GRAUINSTRUCAO_Ensino médio completo",
    "GRAUINSTRUCAO_Ensino médio incompleto",
    "GRAUINSTRUCAO_Mestrado completo",
    "GRAUINSTRUCAO_Pós Grad. completo"
)

# Starting column of the block being rewritten: CR = column 96
$startCol = 96

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $startCol + $i).Value = $newHeaders[$i]
}

# Old block ran through DN1 (column 118); clear the two now-unused trailing
# columns (DM1, DN1) so the used range shrinks back to A1:DL1. Clear() (not
# just blanking the value) drops the cell entirely, including its style, so
# the sheet's used range/dimension shrinks correctly.
$oldEndCol = 118
$clearFrom = $startCol + $newHeaders.Length
for ($c = $clearFrom; $c -le $oldEndCol; $c++) {
    $ws.Cells.Item(1, $c).Clear()
}
